$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target table (header row stays the same; data rows 2-10 are the
# 3x3 cross of Sending/Target cluster in {FAPs, sCs, ECs} for Myoc -> Fzd3).
$rows = @(
  @{ row=2;  A="FAPs"; D="FAPs"; E=2; F=0.6666666666666666; G=0.714093; H=2.142279; I=0.001553869939232348; J=0.001553869939232348; K=2; L=0.6666666666666666; M=0.129292; N=0.387876; O=0.06615700391713267; P=0.06615700391713268; Q=0.09232651215599999; R=0.8309386094039999; S=0.0001027993796565091; T=0.0001027993796565092 },
  @{ row=3;  A="FAPs"; D="sCs";  E=2; F=0.6666666666666666; G=0.714093; H=2.142279; I=0.001553869939232348; J=0.001553869939232348; K=3; L=1;                  M=0.4307096666666667; N=1.292129;  O=0.220388431649395;  P=0.220388431649395;  Q=0.307566757999;      R=2.768100821991;      S=0.0003424549588945578; T=0.0003424549588945578 },
  @{ row=4;  A="FAPs"; D="ECs";  E=2; F=0.6666666666666666; G=0.714093; H=2.142279; I=0.001553869939232348; J=0.001553869939232348; K=3; L=1;                  M=1.394319;           N=4.182957;  O=0.7134545644334723; P=0.7134545644334724; Q=0.995673437667;      R=8.961060939003;      S=0.001108615600681281;  T=0.001108615600681281 },
  @{ row=5;  A="sCs";  D="FAPs"; E=3; F=1;                  G=453.4108886666666; H=1360.232666; I=0.9866243612803347; J=0.9866243612803348; K=2; L=0.6666666666666666; M=0.129292; N=0.387876; O=0.06615700391713267; P=0.06615700391713268; Q=58.62240061749065;  R=527.601605557416;   S=0.06527211173396162; T=0.06527211173396164 },
  @{ row=6;  A="sCs";  D="sCs";  E=3; F=1;                  G=453.4108886666666; H=1360.232666; I=0.9866243612803347; J=0.9866243612803348; K=3; L=1;                  M=0.4307096666666667; N=1.292129; O=0.220388431649395; P=0.220388431649395; Q=195.2884527206571;  R=1757.596074485914;  S=0.217440595609659;   T=0.2174405956096591 },
  @{ row=7;  A="sCs";  D="ECs";  E=3; F=1;                  G=453.4108886666666; H=1360.232666; I=0.9866243612803347; J=0.9866243612803348; K=3; L=1;                  M=1.394319;           N=4.182957;  O=0.7134545644334723; P=0.7134545644334724; Q=632.199416874818;   R=5689.794751873362;  S=0.7039116539367141;  T=0.7039116539367142 },
  @{ row=8;  A="ECs";  D="FAPs"; E=3; F=1;                  G=5.432785666666667; H=16.298357;   I=0.01182176878043295; J=0.01182176878043295; K=2; L=0.6666666666666666; M=0.129292; N=0.387876; O=0.06615700391713267; P=0.06615700391713268; Q=0.7024157244146666; R=6.321741519732;     S=0.000782092803514539; T=0.0007820928035145393 },
  @{ row=9;  A="ECs";  D="sCs";  E=3; F=1;                  G=5.432785666666667; H=16.298357;   I=0.01182176878043295; J=0.01182176878043295; K=3; L=1;                  M=0.4307096666666667; N=1.292129; O=0.220388431649395; P=0.220388431649395; Q=2.339953303561445;  R=21.059579732053;    S=0.002605381080841398; T=0.002605381080841398 },
  @{ row=10; A="ECs";  D="ECs";  E=3; F=1;                  G=5.432785666666667; H=16.298357;   I=0.01182176878043295; J=0.01182176878043295; K=3; L=1;                  M=1.394319;           N=4.182957; O=0.7134545644334723; P=0.7134545644334724; Q=7.575036277961001;  R=68.175326501649;    S=0.008434294896077008; T=0.00843429489607701 }
)

foreach ($item in $rows) {
  $rn = $item.row
  $ws.Range("A$rn").Value = $item.A
  $ws.Range("B$rn").Value = "Myoc"
  $ws.Range("C$rn").Value = "Fzd3"
  $ws.Range("D$rn").Value = $item.D
  $ws.Range("E$rn").Value = $item.E
  $ws.Range("F$rn").Value = $item.F
  $ws.Range("G$rn").Value = $item.G
  $ws.Range("H$rn").Value = $item.H
  $ws.Range("I$rn").Value = $item.I
  $ws.Range("J$rn").Value = $item.J
  $ws.Range("K$rn").Value = $item.K
  $ws.Range("L$rn").Value = $item.L
  $ws.Range("M$rn").Value = $item.M
  $ws.Range("N$rn").Value = $item.N
  $ws.Range("O$rn").Value = $item.O
  $ws.Range("P$rn").Value = $item.P
  $ws.Range("Q$rn").Value = $item.Q
  $ws.Range("R$rn").Value = $item.R
  $ws.Range("S$rn").Value = $item.S
  $ws.Range("T$rn").Value = $item.T
}
